# Balanced reformatted dataset, ptf test.
#
# Adds the missing B-column data points for rows 4 and 5 (rounding out the
# dataset so each labeled row has matching A/B samples) and moves the active
# selection to D7, matching the post-edit sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 0.65234749999999997
$ws.Range("B5").Value = 0.69791668653488104

$ws.Range("D7").Select()
